$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    "D2" = 3.37684444861195
    "E2" = -5.13856388713722
    "F2" = 8.99158073420979
    "G2" = 0
    "D3" = -22.8676613737711
    "E3" = -32.2795665879931
    "F3" = -4.42070723999729
    "G3" = 1
    "D4" = -27.7906514053958
    "E4" = -38.3265017488862
    "F4" = -6.52843663571799
    "G4" = 1
    "D5" = -23.1946597042139
    "E5" = -63.3923204252342
    "F5" = 25.4145741152583
    "D6" = -44.8432635071778
    "E6" = -53.3608177086025
    "F6" = -36.7757597793092
    "D7" = 1.12243934249334
    "E7" = -11.2347177065196
    "F7" = 15.1138482931335
    "G7" = 0
    "D8" = 17.4393956290234
    "E8" = -1.79893411995474
    "F8" = 37.5154343527284
    "G8" = 0
    "D9" = -8.28038947203445
    "E9" = -19.5163408515315
    "F9" = 4.72307014844956
    "G9" = 0
    "D10" = 30.2771796723165
    "E10" = 2.19610836858658
    "F10" = 72.3575384112918
    "D11" = -35.399702251583
    "E11" = -45.6850735016055
    "F11" = -24.4221172343063
    "D12" = 2.583994687556
    "E12" = -2.83512997346274
    "F12" = 15.926019401117
    "D13" = 12.2421704294136
    "E13" = -9.43489729406951
    "F13" = 32.8294924459343
    "G13" = 0
    "D14" = -9.1240236681981
    "E14" = -24.5638658121862
    "F14" = 11.4586564566882
    "G14" = 0
    "D15" = 155.692560143525
    "E15" = -32.4910791518741
    "F15" = 456.821493404389
    "G15" = 0
    "D16" = 144.296975829521
    "E16" = 80.2760983568056
    "F16" = 243.092071745791
    "D17" = -1.85552939885885
    "E17" = -5.30487659806178
    "F17" = 2.34537255745965
    "G17" = 0
    "D18" = 10.6995046902919
    "E18" = -14.42401903237
    "F18" = 28.3404294696575
    "G18" = 0
    "D19" = 26.5230956547787
    "E19" = 3.24128791090334
    "F19" = 83.5406242053416
    "D20" = 34.0914256627225
    "E20" = -3.34258407504167
    "F20" = 102.62496920381
    "G20" = 0
    "D21" = -41.3131527556954
    "E21" = -52.1883779473875
    "F21" = -27.6604294665201
}

foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}
